$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.915.42"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "1.859.36"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("D4").Value = "1.019"
$ws.Range("E4").Value = "  -1.47%  "

$ws.Range("D5").Value = "320.94"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").Value = "1.015"
$ws.Range("E6").Value = "  -1.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4350"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("D8").Value = "0.3789"
$ws.Range("E8").Value = "  -0.24%  "

$ws.Range("D9").Value = "0.07442"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "0.8857"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").Value = "21.65"
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").Value = "1.862.92"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").Value = "6.793"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "5.501"
$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("D15").Value = "0.07146"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("D16").Value = "88.27"
$ws.Range("E16").Value = "  +5.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.020"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009036"
$ws.Range("E18").Value = "  -0.51%  "

$ws.Range("D19").Value = "1.016"
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("D20").Value = "15.45"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "27.913.27"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "5.277"
$ws.Range("E22").Value = "  -0.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.20"
$ws.Range("E23").Value = "  -1.89%  "

$ws.Range("D24").Value = "2.085.51"
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("E25").Value = "  +5.03%  "

$ws.Range("D26").Value = "156.91"
$ws.Range("E26").Value = "  -0.83%  "

$ws.Range("D27").Value = "18.72"
$ws.Range("E27").Value = "  -0.46%  "

$ws.Range("D28").Value = "5.435"
$ws.Range("E28").Value = "  +2.05%  "

$ws.Range("D29").Value = "2.007"
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").Value = "120.39"
$ws.Range("E30").Value = "  +2.35%  "

$ws.Range("D31").Value = "0.09026"
$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").Value = "1.241"
$ws.Range("E32").Value = "  +2.47%  "

$ws.Range("D33").Value = "0.7741"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("D34").Value = "4.576"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("D35").Value = "2.992"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").Value = "1.016"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("D37").Value = "1.144"
$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("D38").Value = "0.01978"
$ws.Range("E38").Value = "  -0.35%  "

$ws.Range("D39").Value = "0.05317"
$ws.Range("E39").Value = "  -0.43%  "

$ws.Range("D40").Value = "2.895"
$ws.Range("E40").Value = "  +2.25%  "

$ws.Range("D41").Value = "0.5203"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.990"
$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("D43").Value = "0.1679"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("D44").Value = "8.747"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("D45").Value = "110.42"
$ws.Range("E45").Value = "  +0.84%  "

$ws.Range("D46").Value = "10.74"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4752"
$ws.Range("E47").Value = "  +1.52%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.714"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.06481"
$ws.Range("E49").Value = "  +0.93%  "

$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "1.016"
$ws.Range("E50").Value = "  -1.63%  "

$ws.Range("D51").Value = "1.877"
$ws.Range("E51").Value = "  +0.55%  "
